$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-28 (Generation 0-26): Fitness column C -> 7295
$ws.Range("C2:C28").Value = 7295

# Rows 29-252 (Generation 27-250): Fitness column C -> 7293
$ws.Range("C29:C252").Value = 7293
